# The deck ships two theme parts:
#   ppt/theme/theme1.xml  ("Office Theme" / clrScheme "Office")   -> used by the Notes Master
#   ppt/theme/theme2.xml  ("Integral" / clrScheme "Red Violet")   -> used by the Slide Master
#
# The target revision swaps the two themes' content: the Slide Master (and
# with it, the single theme color scheme PowerPoint's object model exposes)
# should end up using the default "Office" palette that currently lives in
# theme1.xml, while the "Integral"/"Red Violet" palette moves the other way.
#
# PowerPoint's automation surface only exposes one live ThemeColorScheme for
# the presentation (SlideMaster.Theme / NotesMaster.Theme / Designs(1) all
# resolve to it), so we drive that single object to the "Office" palette --
# dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink, in that fixed Item() order.

function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Target palette = the "Office" color scheme (hex -> RGB() order matches
# Item(1..12): dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4,
# accent5, accent6, hlink, folHlink).
$officePalette = @(
    (RGBVal 0x00 0x00 0x00),   # dk1      000000
    (RGBVal 0xFF 0xFF 0xFF),   # lt1      FFFFFF
    (RGBVal 0x44 0x54 0x6A),   # dk2      44546A
    (RGBVal 0xE7 0xE6 0xE6),   # lt2      E7E6E6
    (RGBVal 0x5B 0x9B 0xD5),   # accent1  5B9BD5
    (RGBVal 0xED 0x7D 0x31),   # accent2  ED7D31
    (RGBVal 0xA5 0xA5 0xA5),   # accent3  A5A5A5
    (RGBVal 0xFF 0xC0 0x00),   # accent4  FFC000
    (RGBVal 0x44 0x72 0xC4),   # accent5  4472C4
    (RGBVal 0x70 0xAD 0x47),   # accent6  70AD47
    (RGBVal 0x05 0x63 0xC1),   # hlink    0563C1
    (RGBVal 0x95 0x4F 0x72)    # folHlink 954F72
)

for ($i = 1; $i -le $officePalette.Length; $i++) {
    $colors.Item($i).RGB = $officePalette[$i - 1]
}
